# Correccion 05 - Finalize OK hasta Cordova
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 3 (header-like row): E3 = "cohecha"
$ws.Range("E3").Value = "cohecha"

# Rows 4-15: fill columns E and F with "ok " (trailing space)
for ($r = 4; $r -le 15; $r++) {
    $ws.Cells.Item($r, 5).Value = "ok "   # column E
    $ws.Cells.Item($r, 6).Value = "ok "   # column F
}

# Row 9 is special: column E should be "ok-" instead of "ok "
$ws.Range("E9").Value = "ok-"

# F3 = "Cordova" (set after "ok-" so shared-string order matches)
$ws.Range("F3").Value = "Cordova"

# Update the active selection to H8, matching the saved selection in the file
$ws.Range("H8").Select()

$wb.Save()
